$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (col D becomes "join" header / boolean FALSE for data rows)
$ws.Range("D1").Value = "join"

$ws.Range("A2").Value = "/mnt/stepanauskas_nfs/julia/testfragrecruitment/hmp/illumina_mgs/SRS019030/SRS019030.denovo_duplicates_marked.trimmed.1.fastq.gz"
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = "illumina"
$ws.Range("D2").Value = $false

$ws.Range("A3").Value = "/mnt/stepanauskas_nfs/julia/testfragrecruitment/hmp/illumina_mgs/SRS022524/SRS022524.denovo_duplicates_marked.trimmed.1.fastq.gz"
$ws.Range("B3").Value = "None"
$ws.Range("C3").Value = "illumina"
$ws.Range("D3").Value = $false

$ws.Range("A4").Value = "/mnt/stepanauskas_nfs/julia/testfragrecruitment/hmp/illumina_mgs/SRS078197/SRS078197.denovo_duplicates_marked.trimmed.1.fastq.gz"
$ws.Range("B4").Value = "None"
$ws.Range("C4").Value = "illumina"
$ws.Range("D4").Value = $false

$ws.Range("A5").Value = "/mnt/stepanauskas_nfs/julia/testfragrecruitment/hmp/pyro_mgs/SRS019030_454/SRS019030_454.fastq.gz"
$ws.Range("B5").Value = "None"
$ws.Range("C5").Value = "pyro"
$ws.Range("D5").Value = $false

$ws.Range("A6").Value = "/mnt/stepanauskas_nfs/julia/testfragrecruitment/hmp/pyro_mgs/SRS022524_454/SRS022524_45.fastq.gz"
$ws.Range("B6").Value = "None"
$ws.Range("C6").Value = "pyro"
$ws.Range("D6").Value = $false

$ws.Range("A7").Value = "/mnt/stepanauskas_nfs/julia/testfragrecruitment/hmp/pyro_mgs/SRS078197_454/SRS078197_454.fastq.gz"
$ws.Range("B7").Value = "None"
$ws.Range("C7").Value = "pyro"
$ws.Range("D7").Value = $false

# Column C width shrinks (no longer needs to fit the long "illumina"/"pyro" + old long text)
# NOTE: Excel's ColumnWidth (chars) differs from the raw OOXML `width` attribute
# by a fixed +5/6 offset for this font/DPI; 41/3 round-trips to exactly 14.5.
$ws.Columns.Item(3).ColumnWidth = 41/3

# Selection moves to D7 (last data row) instead of D8
$ws.Range("D7").Select()

# Window position change
$excel.ActiveWindow.Left = 1680
$excel.ActiveWindow.Top = 0
